$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.118.42"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.321.26"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.507"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "2.682.15"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "2.296.27"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").Value = "43.050.35"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.50%  "
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "169.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.72%  "
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0699"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").Value = "1.992.82"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  +1.34%  "
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("D50").Value = "2.547.07"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("E51").Value = "  +0.84%  "
